$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Table cell "Zach "/"Overy" (split across two runs by a spell-check
#    proofErr wrapper) becomes a single run "Zach Overy" with no
#    proofErr markup.
# ---------------------------------------------------------------------
$rngZach = $d.Content
$null = $rngZach.Find.Execute("Zach Overy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$zachXml = '<?xml version="1.0" encoding="utf-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00ED4CD5" w:rsidRDefault="00ED4CD5" w:rsidP="00ED4CD5"><w:r><w:t>Zach Overy</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $rngZach.InsertXML($zachXml)

# ---------------------------------------------------------------------
# 2. Append a new sentence to the closing paragraph, as its own run
#    (same sz/szCs formatting as the run before it), inserted before
#    the "_GoBack" bookmark that closes the paragraph.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range
$insStart = $bmRange.Start

$newSentence = ". Everyone will be finished coding their individual portions of the project by the middle of November. John will then test everybody" + [char]0x2019 + "s code and we will correct any errors in our code."
$null = $bmRange.InsertBefore($newSentence)

# Force the newly typed text to live in its own run (rather than being
# silently coalesced back into the preceding run) by toggling a
# character-formatting property on it.
$insEnd = $insStart + $newSentence.Length
$insRange = $d.Range($insStart, $insEnd)
$insRange.Bold = $true
$insRange.Bold = $false
